# Auto-generated edit script: applies scheduled-runner market-data refresh
# to the Leve profit calculation sheets (columns H-N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ALC sheet, row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 2902.547  # H17: 2805.3455 -> 2902.547
$ws.Cells.Item(17, 9).Value = 5250  # I17: 2725 -> 5250
$ws.Cells.Item(17, 10).Value = 2857.4038  # J17: 2808.3774 -> 2857.4038
$ws.Cells.Item(17, 11).Value = 15750  # K17: 8175 -> 15750
$ws.Cells.Item(17, 12).Value = 8572.2114  # L17: 8425.1322 -> 8572.2114
$ws.Cells.Item(17, 13).Value = -15582  # M17: -8007 -> -15582
$ws.Cells.Item(17, 14).Value = -8908.2114  # N17: -8761.1322 -> -8908.2114

# ALC sheet, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3658.4138  # H40: 3629.7666 -> 3658.4138
$ws.Cells.Item(40, 9).Value = 2839.8  # I40: 2833 -> 2839.8
$ws.Cells.Item(40, 11).Value = 2839.8  # K40: 2833 -> 2839.8
$ws.Cells.Item(40, 13).Value = -2664.8  # M40: -2658 -> -2664.8

# ALC sheet, row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 140.125  # H42: 140.25 -> 140.125
$ws.Cells.Item(42, 10).Value = 148.6  # J42: 148.8 -> 148.6
$ws.Cells.Item(42, 12).Value = 445.8  # L42: 446.4 -> 445.8
$ws.Cells.Item(42, 14).Value = -905.8  # N42: -906.4000000000001 -> -905.8

# ALC sheet, row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 464.14285  # H55: 383.8 -> 464.14285
$ws.Cells.Item(55, 9).Value = 354.4  # I55: 304.77777 -> 354.4
$ws.Cells.Item(55, 10).Value = 738.5  # J55: 502.33334 -> 738.5
$ws.Cells.Item(55, 11).Value = 354.4  # K55: 304.77777 -> 354.4
$ws.Cells.Item(55, 12).Value = 738.5  # L55: 502.33334 -> 738.5
$ws.Cells.Item(55, 13).Value = -140.4  # M55: -90.77776999999998 -> -140.4
$ws.Cells.Item(55, 14).Value = -1166.5  # N55: -930.33334 -> -1166.5

# ALC sheet, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 1544.1428  # H107: 1500.9546 -> 1544.1428
$ws.Cells.Item(107, 9).Value = 1765.6471  # I107: 1626.8422 -> 1765.6471
$ws.Cells.Item(107, 10).Value = 602.75  # J107: 703.6667 -> 602.75
$ws.Cells.Item(107, 11).Value = 1765.6471  # K107: 1626.8422 -> 1765.6471
$ws.Cells.Item(107, 12).Value = 602.75  # L107: 703.6667 -> 602.75
$ws.Cells.Item(107, 13).Value = 154.3529000000001  # M107: 293.1578 -> 154.3529000000001
$ws.Cells.Item(107, 14).Value = -4442.75  # N107: -4543.6667 -> -4442.75

# ALC sheet, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2646862.2  # H137: 1985293.1 -> 2646862.2
$ws.Cells.Item(137, 9).Value = 896.1111  # I137: 760.3125 -> 896.1111
$ws.Cells.Item(137, 11).Value = 2688.3333  # K137: 2280.9375 -> 2688.3333
$ws.Cells.Item(137, 13).Value = -138.3332999999998  # M137: 269.0625 -> -138.3332999999998

# ALC sheet, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2499.85  # H138: 2461.3625 -> 2499.85
$ws.Cells.Item(138, 10).Value = 2764.8809  # J138: 2749.9866 -> 2764.8809
$ws.Cells.Item(138, 12).Value = 8294.6427  # L138: 8249.959800000001 -> 8294.6427
$ws.Cells.Item(138, 14).Value = -18574.6427  # N138: -18529.9598 -> -18574.6427

# ALC sheet, row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 416666.34  # H140: 319999.75 -> 416666.34
$ws.Cells.Item(140, 10).Value = 416666.34  # J140: 319999.75 -> 416666.34
$ws.Cells.Item(140, 12).Value = 416666.34  # L140: 319999.75 -> 416666.34
$ws.Cells.Item(140, 14).Value = -427026.34  # N140: -330359.75 -> -427026.34

# ARM sheet, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2154.6086  # H2: 2154.6956 -> 2154.6086
$ws.Cells.Item(2, 9).Value = 1772.1428  # I2: 1772.2858 -> 1772.1428
$ws.Cells.Item(2, 11).Value = 1772.1428  # K2: 1772.2858 -> 1772.1428
$ws.Cells.Item(2, 13).Value = -1659.1428  # M2: -1659.2858 -> -1659.1428

# ARM sheet, row 43
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 26000  # H43: 22749.75 -> 26000
$ws.Cells.Item(43, 9).Value = 0  # I43: 12999 -> 0
$ws.Cells.Item(43, 11).Value = 0  # K43: 12999 -> 0
$ws.Cells.Item(43, 13).ClearContents()  # M43: was -12686, now blank

# ARM sheet, row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 15144  # H76: 0 -> 15144
$ws.Cells.Item(76, 10).Value = 15144  # J76: 0 -> 15144
$ws.Cells.Item(76, 12).Value = 15144  # L76: 0 -> 15144
$ws.Cells.Item(76, 14).Value = -15820  # N76: None -> -15820

# ARM sheet, row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 15144  # H79: 0 -> 15144
$ws.Cells.Item(79, 10).Value = 15144  # J79: 0 -> 15144
$ws.Cells.Item(79, 12).Value = 15144  # L79: 0 -> 15144
$ws.Cells.Item(79, 14).Value = -17484  # N79: None -> -17484

# ARM sheet, row 108
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(108, 8).Value = 89291  # H108: 93175.39999999999 -> 89291
$ws.Cells.Item(108, 10).Value = 89291  # J108: 93175.39999999999 -> 89291
$ws.Cells.Item(108, 12).Value = 89291  # L108: 93175.39999999999 -> 89291
$ws.Cells.Item(108, 14).Value = -96971  # N108: -100855.4 -> -96971

# ARM sheet, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2154.6086  # H116: 2154.6956 -> 2154.6086
$ws.Cells.Item(116, 9).Value = 1772.1428  # I116: 1772.2858 -> 1772.1428
$ws.Cells.Item(116, 11).Value = 1772.1428  # K116: 1772.2858 -> 1772.1428
$ws.Cells.Item(116, 13).Value = 521.8571999999999  # M116: 521.7141999999999 -> 521.8571999999999

# ARM sheet, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4193.0645  # H132: 3731.9487 -> 4193.0645
$ws.Cells.Item(132, 9).Value = 4040.037  # I132: 3561.2 -> 4040.037
$ws.Cells.Item(132, 11).Value = 12120.111  # K132: 10683.6 -> 12120.111
$ws.Cells.Item(132, 13).Value = -9590.110999999999  # M132: -8153.599999999999 -> -9590.110999999999

# BSM sheet, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2154.6086  # H3: 2154.6956 -> 2154.6086
$ws.Cells.Item(3, 9).Value = 1772.1428  # I3: 1772.2858 -> 1772.1428
$ws.Cells.Item(3, 11).Value = 1772.1428  # K3: 1772.2858 -> 1772.1428
$ws.Cells.Item(3, 13).Value = -1658.1428  # M3: -1658.2858 -> -1658.1428

# BSM sheet, row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1236.5312  # H86: 1304.2667 -> 1236.5312
$ws.Cells.Item(86, 9).Value = 1202.2632  # I86: 1207.7894 -> 1202.2632
$ws.Cells.Item(86, 10).Value = 1286.6154  # J86: 1470.909 -> 1286.6154
$ws.Cells.Item(86, 11).Value = 1202.2632  # K86: 1207.7894 -> 1202.2632
$ws.Cells.Item(86, 12).Value = 1286.6154  # L86: 1470.909 -> 1286.6154
$ws.Cells.Item(86, 13).Value = -79.2632000000001  # M86: -84.78939999999989 -> -79.2632000000001
$ws.Cells.Item(86, 14).Value = -3532.6154  # N86: -3716.909 -> -3532.6154

# BSM sheet, row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 1236.5312  # H89: 1304.2667 -> 1236.5312
$ws.Cells.Item(89, 9).Value = 1202.2632  # I89: 1207.7894 -> 1202.2632
$ws.Cells.Item(89, 10).Value = 1286.6154  # J89: 1470.909 -> 1286.6154
$ws.Cells.Item(89, 11).Value = 6011.316000000001  # K89: 6038.946999999999 -> 6011.316000000001
$ws.Cells.Item(89, 12).Value = 6433.076999999999  # L89: 7354.545 -> 6433.076999999999
$ws.Cells.Item(89, 13).Value = -395.3160000000007  # M89: -422.9469999999992 -> -395.3160000000007
$ws.Cells.Item(89, 14).Value = -17665.077  # N89: -18586.545 -> -17665.077

# BSM sheet, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1929.3334  # H94: 3111.25 -> 1929.3334
$ws.Cells.Item(94, 9).Value = 1473.5714  # I94: 1763 -> 1473.5714
$ws.Cells.Item(94, 10).Value = 2328.125  # J94: 5358.3335 -> 2328.125
$ws.Cells.Item(94, 11).Value = 1473.5714  # K94: 1763 -> 1473.5714
$ws.Cells.Item(94, 12).Value = 2328.125  # L94: 5358.3335 -> 2328.125
$ws.Cells.Item(94, 13).Value = -1022.5714  # M94: -1312 -> -1022.5714
$ws.Cells.Item(94, 14).Value = -3230.125  # N94: -6260.3335 -> -3230.125

# CRP sheet, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3924.5  # H31: 4000.7693 -> 3924.5
$ws.Cells.Item(31, 9).Value = 1439.8889  # I31: 1501.125 -> 1439.8889
$ws.Cells.Item(31, 11).Value = 1439.8889  # K31: 1501.125 -> 1439.8889
$ws.Cells.Item(31, 13).Value = -1144.8889  # M31: -1206.125 -> -1144.8889

# CRP sheet, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3924.5  # H34: 4000.7693 -> 3924.5
$ws.Cells.Item(34, 9).Value = 1439.8889  # I34: 1501.125 -> 1439.8889
$ws.Cells.Item(34, 11).Value = 1439.8889  # K34: 1501.125 -> 1439.8889
$ws.Cells.Item(34, 13).Value = -1237.8889  # M34: -1299.125 -> -1237.8889

# CRP sheet, row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(44, 8).Value = 26000  # H44: 0 -> 26000
$ws.Cells.Item(44, 10).Value = 26000  # J44: 0 -> 26000
$ws.Cells.Item(44, 12).Value = 26000  # L44: 0 -> 26000
$ws.Cells.Item(44, 14).Value = -26884  # N44: None -> -26884

# CRP sheet, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2793.8696  # H58: 2554.0195 -> 2793.8696
$ws.Cells.Item(58, 9).Value = 2499.75  # I58: 2260.6 -> 2499.75
$ws.Cells.Item(58, 11).Value = 2499.75  # K58: 2260.6 -> 2499.75
$ws.Cells.Item(58, 13).Value = -2296.75  # M58: -2057.6 -> -2296.75

# CRP sheet, row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(70, 8).Value = 36999.668  # H70: 37000 -> 36999.668
$ws.Cells.Item(70, 10).Value = 36999.668  # J70: 37000 -> 36999.668
$ws.Cells.Item(70, 12).Value = 36999.668  # L70: 37000 -> 36999.668
$ws.Cells.Item(70, 14).Value = -37629.668  # N70: -37630 -> -37629.668

# CRP sheet, row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(73, 8).Value = 36999.668  # H73: 37000 -> 36999.668
$ws.Cells.Item(73, 10).Value = 36999.668  # J73: 37000 -> 36999.668
$ws.Cells.Item(73, 12).Value = 36999.668  # L73: 37000 -> 36999.668
$ws.Cells.Item(73, 14).Value = -39183.668  # N73: -39184 -> -39183.668

# CRP sheet, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2793.8696  # H136: 2554.0195 -> 2793.8696
$ws.Cells.Item(136, 9).Value = 2499.75  # I136: 2260.6 -> 2499.75
$ws.Cells.Item(136, 11).Value = 7499.25  # K136: 6781.799999999999 -> 7499.25
$ws.Cells.Item(136, 13).Value = -4949.25  # M136: -4231.799999999999 -> -4949.25

# CUL sheet, row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 833408.5600000001  # H12: 1111220 -> 833408.5600000001
$ws.Cells.Item(12, 9).Value = 27.5  # I12: 0 -> 27.5
$ws.Cells.Item(12, 10).Value = 1000084.8  # J12: 1111220 -> 1000084.8
$ws.Cells.Item(12, 11).Value = 82.5  # K12: 0 -> 82.5
$ws.Cells.Item(12, 12).Value = 3000254.4  # L12: 3333660 -> 3000254.4
$ws.Cells.Item(12, 13).Value = 90.5  # M12: None -> 90.5
$ws.Cells.Item(12, 14).Value = -3000600.4  # N12: -3334006 -> -3000600.4

# CUL sheet, row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 1643  # H14: 2095 -> 1643
$ws.Cells.Item(14, 9).Value = 1643  # I14: 2095 -> 1643
$ws.Cells.Item(14, 11).Value = 4929  # K14: 6285 -> 4929
$ws.Cells.Item(14, 13).Value = -4756  # M14: -6112 -> -4756

# CUL sheet, row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 3742.2666  # H18: 3867.9285 -> 3742.2666
$ws.Cells.Item(18, 9).Value = 2792.889  # I18: 2794.7778 -> 2792.889
$ws.Cells.Item(18, 10).Value = 5166.3335  # J18: 5799.6 -> 5166.3335
$ws.Cells.Item(18, 11).Value = 8378.667000000001  # K18: 8384.3334 -> 8378.667000000001
$ws.Cells.Item(18, 12).Value = 15499.0005  # L18: 17398.8 -> 15499.0005
$ws.Cells.Item(18, 13).Value = -8209.667000000001  # M18: -8215.3334 -> -8209.667000000001
$ws.Cells.Item(18, 14).Value = -15837.0005  # N18: -17736.8 -> -15837.0005

# CUL sheet, row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 46.65  # H38: 51.055557 -> 46.65
$ws.Cells.Item(38, 9).Value = 44.25  # I38: 50.9 -> 44.25
$ws.Cells.Item(38, 10).Value = 50.25  # J38: 51.25 -> 50.25
$ws.Cells.Item(38, 11).Value = 132.75  # K38: 152.7 -> 132.75
$ws.Cells.Item(38, 12).Value = 150.75  # L38: 153.75 -> 150.75
$ws.Cells.Item(38, 13).Value = 214.25  # M38: 194.3 -> 214.25
$ws.Cells.Item(38, 14).Value = -844.75  # N38: -847.75 -> -844.75

# CUL sheet, row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 3551.4285  # H39: 3481.2666 -> 3551.4285
$ws.Cells.Item(39, 10).Value = 3551.4285  # J39: 3481.2666 -> 3551.4285
$ws.Cells.Item(39, 12).Value = 10654.2855  # L39: 10443.7998 -> 10654.2855
$ws.Cells.Item(39, 14).Value = -11242.2855  # N39: -11031.7998 -> -11242.2855

# CUL sheet, row 124
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value = 1999.6666  # H124: 2000 -> 1999.6666
$ws.Cells.Item(124, 9).Value = 1999.6666  # I124: 2000 -> 1999.6666
$ws.Cells.Item(124, 11).Value = 5998.9998  # K124: 6000 -> 5998.9998
$ws.Cells.Item(124, 13).Value = -1088.9998  # M124: -1090 -> -1088.9998

# CUL sheet, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1679.4147  # H131: 1724.1052 -> 1679.4147
$ws.Cells.Item(131, 9).Value = 1367.3334  # I131: 1451.8 -> 1367.3334
$ws.Cells.Item(131, 10).Value = 1732.9143  # J131: 1765.3636 -> 1732.9143
$ws.Cells.Item(131, 11).Value = 4102.0002  # K131: 4355.4 -> 4102.0002
$ws.Cells.Item(131, 12).Value = 5198.742899999999  # L131: 5296.0908 -> 5198.742899999999
$ws.Cells.Item(131, 13).Value = 937.9997999999996  # M131: 684.6000000000004 -> 937.9997999999996
$ws.Cells.Item(131, 14).Value = -15278.7429  # N131: -15376.0908 -> -15278.7429

# GSM sheet, row 49
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(49, 8).Value = 30000  # H49: 0 -> 30000
$ws.Cells.Item(49, 10).Value = 30000  # J49: 0 -> 30000
$ws.Cells.Item(49, 12).Value = 30000  # L49: 0 -> 30000
$ws.Cells.Item(49, 14).Value = -30368  # N49: None -> -30368

# GSM sheet, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 688.44446  # H107: 831.7143 -> 688.44446
$ws.Cells.Item(107, 9).Value = 688.8  # I107: 817.25 -> 688.8
$ws.Cells.Item(107, 10).Value = 688  # J107: 851 -> 688
$ws.Cells.Item(107, 11).Value = 688.8  # K107: 817.25 -> 688.8
$ws.Cells.Item(107, 12).Value = 688  # L107: 851 -> 688
$ws.Cells.Item(107, 13).Value = 1231.2  # M107: 1102.75 -> 1231.2
$ws.Cells.Item(107, 14).Value = -4528  # N107: -4691 -> -4528

# GSM sheet, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1861.8889  # H122: 1960.2727 -> 1861.8889
$ws.Cells.Item(122, 9).Value = 1625.75  # I122: 1884.8334 -> 1625.75
$ws.Cells.Item(122, 11).Value = 4877.25  # K122: 5654.5002 -> 4877.25
$ws.Cells.Item(122, 13).Value = -2427.25  # M122: -3204.5002 -> -2427.25

# LTW sheet, row 42
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 17303  # H42: 18166.666 -> 17303
$ws.Cells.Item(42, 9).Value = 15508.333  # I42: 17250 -> 15508.333
$ws.Cells.Item(42, 10).Value = 19995  # J42: 20000 -> 19995
$ws.Cells.Item(42, 11).Value = 15508.333  # K42: 17250 -> 15508.333
$ws.Cells.Item(42, 12).Value = 19995  # L42: 20000 -> 19995
$ws.Cells.Item(42, 13).Value = -14945.333  # M42: -16687 -> -14945.333
$ws.Cells.Item(42, 14).Value = -21121  # N42: -21126 -> -21121

# LTW sheet, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 10030.235  # H46: 10594.625 -> 10030.235
$ws.Cells.Item(46, 9).Value = 1000  # I46: 0 -> 1000
$ws.Cells.Item(46, 11).Value = 1000  # K46: 0 -> 1000
$ws.Cells.Item(46, 13).Value = -812  # M46: None -> -812

# LTW sheet, row 49
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(49, 8).Value = 17303  # H49: 18166.666 -> 17303
$ws.Cells.Item(49, 9).Value = 15508.333  # I49: 17250 -> 15508.333
$ws.Cells.Item(49, 10).Value = 19995  # J49: 20000 -> 19995
$ws.Cells.Item(49, 11).Value = 15508.333  # K49: 17250 -> 15508.333
$ws.Cells.Item(49, 12).Value = 19995  # L49: 20000 -> 19995
$ws.Cells.Item(49, 13).Value = -15361.333  # M49: -17103 -> -15361.333
$ws.Cells.Item(49, 14).Value = -20289  # N49: -20294 -> -20289

# LTW sheet, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1459.4762  # H61: 1498 -> 1459.4762
$ws.Cells.Item(61, 9).Value = 1366.7646  # I61: 1409.125 -> 1366.7646
$ws.Cells.Item(61, 11).Value = 1366.7646  # K61: 1409.125 -> 1366.7646
$ws.Cells.Item(61, 13).Value = -1164.7646  # M61: -1207.125 -> -1164.7646

# LTW sheet, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3353.2666  # H93: 3438.4614 -> 3353.2666
$ws.Cells.Item(93, 9).Value = 2943.4285  # I93: 2950.8333 -> 2943.4285
$ws.Cells.Item(93, 10).Value = 3711.875  # J93: 3856.4285 -> 3711.875
$ws.Cells.Item(93, 11).Value = 2943.4285  # K93: 2950.8333 -> 2943.4285
$ws.Cells.Item(93, 12).Value = 3711.875  # L93: 3856.4285 -> 3711.875
$ws.Cells.Item(93, 13).Value = -1695.4285  # M93: -1702.8333 -> -1695.4285
$ws.Cells.Item(93, 14).Value = -6207.875  # N93: -6352.4285 -> -6207.875

# LTW sheet, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 0  # H100: 2500 -> 0
$ws.Cells.Item(100, 10).Value = 0  # J100: 2500 -> 0
$ws.Cells.Item(100, 12).Value = 0  # L100: 2500 -> 0
$ws.Cells.Item(100, 14).ClearContents()  # N100: was -3582, now blank

# LTW sheet, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1459.4762  # H113: 1498 -> 1459.4762
$ws.Cells.Item(113, 9).Value = 1366.7646  # I113: 1409.125 -> 1366.7646
$ws.Cells.Item(113, 11).Value = 1366.7646  # K113: 1409.125 -> 1366.7646
$ws.Cells.Item(113, 13).Value = 803.2354  # M113: 760.875 -> 803.2354

# LTW sheet, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 13608.923  # H136: 12735.429 -> 13608.923
$ws.Cells.Item(136, 9).Value = 5319.6665  # I136: 4756.857 -> 5319.6665
$ws.Cells.Item(136, 11).Value = 15958.9995  # K136: 14270.571 -> 15958.9995
$ws.Cells.Item(136, 13).Value = -13408.9995  # M136: -11720.571 -> -13408.9995
